# Update the handoff/handback timestamps for the zh-cn and de-de report sheets.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-28 09:15:08"
$wsZhCn.Range("G3").Value = "2016-01-28 09:15:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-28 09:15:21"
$wsDeDe.Range("G3").Value = "2016-01-28 09:16:15"
